# IMU fusion decision matrix - add "Complementary Filter" solution column,
# rename "Madgwick's Filter" -> "Madgwick Filter", fix DMP/Ease-of-use score,
# and refresh the table look (bold headers, zoom, borders).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column D ("Complementary Filter") - shifts old D:G -> E:H
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).Insert()

# Rename the Madgwick column header text (now in column G after the insert)
$ws.Range("G3").Value2 = "Madgwick Filter"

# New "Complementary Filter" header + scores
$ws.Range("D3").Value2 = "Complementary Filter"
$ws.Range("D4").Value2 = 0
$ws.Range("D5").Value2 = -1
$ws.Range("D6").Value2 = -1
$ws.Range("D7").Value2 = -1
$ws.Range("D8").Value2 = 0

# Fix up the DMP column (now H) - "Ease of use" score changed 0 -> 1
$ws.Range("H4").Value2 = 1

# ---------------------------------------------------------------------------
# 2. Totals row - add formula for new column D, keep others consistent
# ---------------------------------------------------------------------------
$ws.Range("D9").FormulaArray = "=SUM(C4:C8*D4:D8)"
$ws.Range("E9").FormulaArray = "=SUM(C4:C8*E4:E8)"
$ws.Range("F9").FormulaArray = "=SUM(C4:C8*F4:F8)"
$ws.Range("G9").FormulaArray = "=SUM(C4:C8*G4:G8)"
$ws.Range("H9").FormulaArray = "=SUM(H4:H8*C4:C8)"

# ---------------------------------------------------------------------------
# 3. Merged header cell now spans D2:H2
# ---------------------------------------------------------------------------
$ws.Range("D2:H2").Merge()

# ---------------------------------------------------------------------------
# 4. Column widths / row heights
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 23.5703125
$ws.Columns.Item(4).ColumnWidth = 14.5703125
$ws.Range("E1:H1").ColumnWidth = 11.7109375
$ws.Rows.Item(3).RowHeight = 45

# ---------------------------------------------------------------------------
# 5. Header formatting: bold font across the title rows + totals label
# ---------------------------------------------------------------------------
$ws.Range("B2:H3").Font.Bold = $true
$ws.Range("B9:C9").Font.Bold = $true

$ws.Range("B2:H3").HorizontalAlignment = -4108   # xlCenter
$ws.Range("B2:H3").VerticalAlignment = -4108     # xlCenter
$ws.Range("D3:H3").WrapText = $true

$ws.Range("B9:C9").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 6. Borders - full grid across the table, outside box + inside gridlines
# ---------------------------------------------------------------------------
$tbl = $ws.Range("B2:H9")
$tbl.Borders.Item(7).LineStyle = 1    # xlEdgeLeft
$tbl.Borders.Item(7).Weight = 2
$tbl.Borders.Item(8).LineStyle = 1    # xlEdgeTop
$tbl.Borders.Item(8).Weight = 2
$tbl.Borders.Item(9).LineStyle = 1    # xlEdgeBottom
$tbl.Borders.Item(9).Weight = 2
$tbl.Borders.Item(10).LineStyle = 1   # xlEdgeRight
$tbl.Borders.Item(10).Weight = 2
$tbl.Borders.Item(11).LineStyle = 1   # xlInsideVertical
$tbl.Borders.Item(11).Weight = 2
$tbl.Borders.Item(12).LineStyle = 1   # xlInsideHorizontal
$tbl.Borders.Item(12).Weight = 2

# ---------------------------------------------------------------------------
# 7. Conditional formatting - colour scale now covers D4:H8, plus an extra
#    colour scale rule restricted to E4:H8 (matches the authored diff)
# ---------------------------------------------------------------------------
$ws.Range("D4:H8").FormatConditions.Delete()

$cf1 = $ws.Range("D4:H8").FormatConditions.AddColorScale(2)
$cf1.ColorScaleCriteria.Item(1).Type = 1     # xlConditionValueLowestValue
$cf1.ColorScaleCriteria.Item(1).FormatColor.Color = 16749212
$cf1.ColorScaleCriteria.Item(2).Type = 2     # xlConditionValueHighestValue
$cf1.ColorScaleCriteria.Item(2).FormatColor.Color = 6534524

$cf2 = $ws.Range("E4:H8").FormatConditions.AddColorScale(2)
$cf2.ColorScaleCriteria.Item(1).Type = 1
$cf2.ColorScaleCriteria.Item(1).FormatColor.Color = 16749212
$cf2.ColorScaleCriteria.Item(2).Type = 2
$cf2.ColorScaleCriteria.Item(2).FormatColor.Color = 6534524

# ---------------------------------------------------------------------------
# 8. View settings - zoomed in, selection moved to D2
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 205
$ws.Range("D2").Select()

$wb.Save()
